$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 196, pushing existing rows 196:200 down to 197:201
$ws.Rows("196").Insert()

# Populate the newly inserted row 196 with the new record
$ws.Range("A196").Value = 4
$ws.Range("B196").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C196").Value = 'Los Lagos'
$ws.Range("D196").Value = 44610
$ws.Range("E196").Value = 10
$ws.Range("F196").Value = 'Fruta'
$ws.Range("G196").Value = 100109
$ws.Range("H196").Value = 'Uva'
$ws.Range("I196").Value = 100109001
$ws.Range("J196").Value = 'Uva'
$ws.Range("K196").Value = 'Red Globe'
$ws.Range("L196").Value = 'Primera'
$ws.Range("M196").Value = 300
$ws.Range("N196").Value = 17000
$ws.Range("O196").Value = 17000
$ws.Range("P196").Value = 17000
$ws.Range("Q196").Value = '$/caja 20 kilos'
$ws.Range("R196").Value = "Región de O'Higgins"
$ws.Range("S196").Value = 850
$ws.Range("T196").Value = 20
